# Apply cell updates per the crypto-price refresh diff.
# Leading apostrophe forces text interpretation (source cells are
# inline strings, even for numeric-looking content like "1.01").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "26.620.64"
$ws.Range("E2").Value = "'" + "  -0.07%  "
$ws.Range("D3").Value = "'" + "1.643.15"
$ws.Range("E3").Value = "'" + "  +0.66%  "
$ws.Range("D4").Value = "'" + "1.01"
$ws.Range("E4").Value = "'" + "  +0.25%  "
$ws.Range("D5").Value = "'" + "215.94"
$ws.Range("E5").Value = "'" + "  +1.50%  "
$ws.Range("E6").Value = "'" + "  +1.12%  "
$ws.Range("E7").Value = "'" + "  +0.22%  "
$ws.Range("E9").Value = "'" + "  +0.83%  "
$ws.Range("E10").Value = "'" + "  +0.48%  "
$ws.Range("D11").Value = "'" + "0.0844"
$ws.Range("E11").Value = "'" + "  +0.04%  "
$ws.Range("D12").Value = "'" + "1.872.67"
$ws.Range("E12").Value = "'" + "  +0.67%  "
$ws.Range("B13").Value = "'" + "WrappedEther"
$ws.Range("C13").Value = "'" + "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'" + "1.678.22"
$ws.Range("E13").Value = "'" + "  +3.52%  "
$ws.Range("B14").Value = "'" + "Polkadot"
$ws.Range("C14").Value = "'" + "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'" + "4.21"
$ws.Range("E14").Value = "'" + "  +3.11%  "
$ws.Range("E15").Value = "'" + "  +1.49%  "
$ws.Range("D16").Value = "'" + "65.91"
$ws.Range("D17").Value = "'" + "26.681.03"
$ws.Range("E17").Value = "'" + "  +0.16%  "
$ws.Range("E18").Value = "'" + "  +1.58%  "
$ws.Range("D19").Value = "'" + "218.37"
$ws.Range("E19").Value = "'" + "  +0.20%  "
$ws.Range("E20").Value = "'" + "  +0.25%  "
$ws.Range("D21").Value = "'" + "4.38"
$ws.Range("E21").Value = "'" + "  +2.36%  "
$ws.Range("E22").Value = "'" + "  +2.22%  "
$ws.Range("E24").Value = "'" + "  +10.05%  "
$ws.Range("D25").Value = "'" + "146.27"
$ws.Range("E25").Value = "'" + "  -1.16%  "
$ws.Range("E26").Value = "'" + "  +0.28%  "
$ws.Range("E27").Value = "'" + "  -0.21%  "
$ws.Range("E28").Value = "'" + "  +3.81%  "
$ws.Range("D29").Value = "'" + "15.85"
$ws.Range("E29").Value = "'" + "  +2.43%  "
$ws.Range("E30").Value = "'" + "  +2.85%  "
$ws.Range("E31").Value = "'" + "  +0.83%  "
$ws.Range("E32").Value = "'" + "  +2.99%  "
$ws.Range("E33").Value = "'" + "  +2.49%  "
$ws.Range("D34").Value = "'" + "1.274.38"
$ws.Range("E34").Value = "'" + "  +5.37%  "
$ws.Range("E36").Value = "'" + "  +6.21%  "
$ws.Range("E38").Value = "'" + "  +5.96%  "
$ws.Range("D39").Value = "'" + "0.830"
$ws.Range("E39").Value = "'" + "  +2.69%  "
$ws.Range("E40").Value = "'" + "  +0.28%  "
$ws.Range("E41").Value = "'" + "  +2.38%  "
$ws.Range("E42").Value = "'" + "  -1.41%  "
$ws.Range("D44").Value = "'" + "1.784.01"
$ws.Range("E44").Value = "'" + "  +0.62%  "
$ws.Range("E45").Value = "'" + "  +0.38%  "
$ws.Range("D46").Value = "'" + "59.76"
$ws.Range("E46").Value = "'" + "  +9.32%  "
$ws.Range("D47").Value = "'" + "1.60"
$ws.Range("E47").Value = "'" + "  +3.83%  "
$ws.Range("E48").Value = "'" + "  +0.82%  "
$ws.Range("D49").Value = "'" + "7.81"
$ws.Range("E49").Value = "'" + "  +3.02%  "
$ws.Range("D50").Value = "'" + "0.0977"
$ws.Range("E50").Value = "'" + "  +4.04%  "
$ws.Range("E51").Value = "'" + "  -0.42%  "
